$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 188; existing rows 188-221 shift down to 189-222.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new data record.
$ws.Cells.Item(188, 1).Value = 9
$ws.Cells.Item(188, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44722
$ws.Cells.Item(188, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = 100112026
$ws.Cells.Item(188, 7).Value = "Haba"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 52
$ws.Cells.Item(188, 11).Value = 18000
$ws.Cells.Item(188, 12).Value = 20000
$ws.Cells.Item(188, 13).Value = 19000
$ws.Cells.Item(188, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(188, 16).Value = 760
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"
